$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 ("time_taken"), matching the style of the
# existing header row (bold, centered, bordered - same as E1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Populate the new time_taken column for each data row (plain/default style).
$ws.Range("F2").Value = "2021-10-05 10:51:38.521469"
$ws.Range("F3").Value = "2021-10-05 10:51:38.521479"
$ws.Range("F4").Value = "2021-10-05 10:51:38.521483"
$ws.Range("F5").Value = "2021-10-05 10:51:38.521485"
$ws.Range("F6").Value = "2021-10-05 10:51:38.521488"
$ws.Range("F7").Value = "2021-10-05 10:51:38.521491"
